$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6142318033764959
$ws.Range("J2").Value = 0.6142318033764957
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.1940157382453333
$ws.Range("R2").Value = 1.746141644208
$ws.Range("S2").Value = 0.00691778299350241
$ws.Range("T2").Value = 0.00691778299350241

# Row 3
$ws.Range("I3").Value = 0.6142318033764959
$ws.Range("J3").Value = 0.6142318033764957
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("S3").Value = 0.5639296366787447
$ws.Range("T3").Value = 0.5639296366787447

# Row 4
$ws.Range("I4").Value = 0.6142318033764959
$ws.Range("J4").Value = 0.6142318033764957
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 1.216755894280666
$ws.Range("R4").Value = 10.950803048526
$ws.Range("S4").Value = 0.04338438370424865
$ws.Range("T4").Value = 0.04338438370424865

# Row 5
$ws.Range("G5").Value = 1.102210333333334
$ws.Range("H5").Value = 3.306631
$ws.Range("I5").Value = 0.3857681966235041
$ws.Range("J5").Value = 0.3857681966235041
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 0.1218515567706667
$ws.Range("R5").Value = 1.096664010936
$ws.Range("S5").Value = 0.004344712623746062
$ws.Range("T5").Value = 0.004344712623746062

# Row 6
$ws.Range("G6").Value = 1.102210333333334
$ws.Range("H6").Value = 3.306631
$ws.Range("I6").Value = 0.3857681966235041
$ws.Range("J6").Value = 0.3857681966235041
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("Q6").Value = 9.933197413530223
$ws.Range("R6").Value = 89.39877672177201
$ws.Range("S6").Value = 0.3541759279936886
$ws.Range("T6").Value = 0.3541759279936885

# Row 7
$ws.Range("G7").Value = 1.102210333333334
$ws.Range("H7").Value = 3.306631
$ws.Range("I7").Value = 0.3857681966235041
$ws.Range("J7").Value = 0.3857681966235041
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 0.7641833660963333
$ws.Range("R7").Value = 6.877650294867
$ws.Range("S7").Value = 0.02724755600606951
$ws.Range("T7").Value = 0.02724755600606951
